# Generate Report for Handback
#
# - Marks the two source files as "Handed back: in sync with en-US" on the
#   Overview sheet as well as the per-locale (zh-cn / de-de) sheets.
# - Populates the "Latest Target File" / "Latest Handback File" columns
#   (F/G) for both rows on the zh-cn and de-de sheets, re-using the same
#   hyperlink targets already used by the Source File Name / Latest
#   Handoff File columns.
# - Stamps the "Latest Handback DateTime" column (H) with the actual
#   handback timestamps.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# hyperlink-style font color (matches the workbook's existing custom
# "HyperLink" cell style: font color FF6495ED, single underline)
$linkColor = 15570276  # BGR long for RGB(0x64,0x95,0xED)

function Style-AsLink($range) {
    $range.Font.Underline = 2      # xlUnderlineStyleSingle
    $range.Font.Color = $linkColor
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
}

# ---------------------------------------------------------------------
# Overview sheet: update Status columns (B, C) for both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

# Row 2 - Latest Target File (F2) / Latest Handback File (G2)
$zhcn.Range("F2").Value = "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0dd42b69bcd9dc5bdd1e521c88bb04cd180c3f6/e2e/35c83f7e-80e6-40ab-987c-2d98ccfc9dae.md", "", "", "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.md") | Out-Null
Style-AsLink $zhcn.Range("F2")

$zhcn.Range("G2").Value = "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.e8b1edb86afeb16dcce65437cc2657b495642639.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9baee7a2ad2d57753a079cd7ccce15ad5d78bf84/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/35c83f7e-80e6-40ab-987c-2d98ccfc9dae.e8b1edb86afeb16dcce65437cc2657b495642639.zh-cn.xlf", "", "", "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.e8b1edb86afeb16dcce65437cc2657b495642639.zh-cn.xlf") | Out-Null
Style-AsLink $zhcn.Range("G2")

# Row 2 - Latest Handback DateTime (H2)
$zhcn.Range("H2").Value = "2016-03-24 11:21:31"

# Row 3 - Latest Target File (F3) / Latest Handback File (G3)
$zhcn.Range("F3").Value = "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/c0dd42b69bcd9dc5bdd1e521c88bb04cd180c3f6/e2e/b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.md", "", "", "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.md") | Out-Null
Style-AsLink $zhcn.Range("F3")

$zhcn.Range("G3").Value = "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.790623e1e94e25f3ddc9e8d04f17f60ae242c55e.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9baee7a2ad2d57753a079cd7ccce15ad5d78bf84/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.790623e1e94e25f3ddc9e8d04f17f60ae242c55e.zh-cn.xlf", "", "", "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.790623e1e94e25f3ddc9e8d04f17f60ae242c55e.zh-cn.xlf") | Out-Null
Style-AsLink $zhcn.Range("G3")

# Row 3 - Latest Handback DateTime (H3)
$zhcn.Range("H3").Value = "2016-03-24 11:21:31"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# Row 2 - Latest Target File (F2) / Latest Handback File (G2)
$dede.Range("F2").Value = "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0dd42b69bcd9dc5bdd1e521c88bb04cd180c3f6/e2e/35c83f7e-80e6-40ab-987c-2d98ccfc9dae.md", "", "", "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.md") | Out-Null
Style-AsLink $dede.Range("F2")

$dede.Range("G2").Value = "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.e8b1edb86afeb16dcce65437cc2657b495642639.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2330c75ccacbe6d644df9c0d646c9d1db1736120/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/35c83f7e-80e6-40ab-987c-2d98ccfc9dae.e8b1edb86afeb16dcce65437cc2657b495642639.de-de.xlf", "", "", "35c83f7e-80e6-40ab-987c-2d98ccfc9dae.e8b1edb86afeb16dcce65437cc2657b495642639.de-de.xlf") | Out-Null
Style-AsLink $dede.Range("G2")

# Row 2 - Latest Handback DateTime (H2)
$dede.Range("H2").Value = "2016-03-24 11:21:46"

# Row 3 - Latest Target File (F3) / Latest Handback File (G3)
$dede.Range("F3").Value = "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/c0dd42b69bcd9dc5bdd1e521c88bb04cd180c3f6/e2e/b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.md", "", "", "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.md") | Out-Null
Style-AsLink $dede.Range("F3")

$dede.Range("G3").Value = "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.790623e1e94e25f3ddc9e8d04f17f60ae242c55e.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2330c75ccacbe6d644df9c0d646c9d1db1736120/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.790623e1e94e25f3ddc9e8d04f17f60ae242c55e.de-de.xlf", "", "", "b6c36069-db4a-4ca0-9f4d-5e5c3a663aeb.790623e1e94e25f3ddc9e8d04f17f60ae242c55e.de-de.xlf") | Out-Null
Style-AsLink $dede.Range("G3")

# Row 3 - Latest Handback DateTime (H3)
$dede.Range("H3").Value = "2016-03-24 11:21:46"

Write-Host "Handback report generated."
